$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.571.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.56%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.025.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.33%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '379.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.03%  '

# Row 7
$ws.Range("E7").Value = '  +1.46%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.78%  '

# Row 11
$ws.Range("E11").Value = '  -0.16%  '

# Row 12
$ws.Range("E12").Value = '  +0.82%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.497.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.36%  '

# Row 15
$ws.Range("E15").Value = '  +1.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.017.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.32%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.983'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -14.38%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.589.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.65%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.64%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.16%  '

# Row 22
$ws.Range("E22").Value = '  +1.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.79%  '

# Row 25
$ws.Range("E25").Value = '  -0.89%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.34%  '

# Row 27
$ws.Range("E27").Value = '  +5.42%  '

# Row 28
$ws.Range("E28").Value = '  +6.60%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.34%  '

# Row 30
$ws.Range("E30").Value = '  -0.05%  '

# Row 31
$ws.Range("E31").Value = '  +1.58%  '

# Row 32
$ws.Range("E32").Value = '  +3.31%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.94%  '

# Row 34
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.40%  '

# Row 35
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '

# Row 36
$ws.Range("E36").Value = '  +5.37%  '

# Row 37
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
$ws.Range("E38").Value = '  +6.77%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.02%  '

# Row 40
$ws.Range("E40").Value = '  +4.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.284'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.17%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '128.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.51%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.68%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.78%  '

# Row 47
$ws.Range("E47").Value = '  +1.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.78%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.027.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.43%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.323.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0320'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.00%  '
